# Edit script: applies the "User Story.docx" diff via Word COM-interop.
#
# Strategy: for the paragraphs whose *internal run structure* changes
# (not just a plain text swap), we rebuild the paragraph's runs with
# InsertXML so we get exact control over run boundaries / proofErr
# markers without the engine's "merge adjacent same-format runs"
# normalization getting in the way. For the two simple sentence-level
# wording tweaks we use a plain Find/Replace.

$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function New-Pkg([string]$innerBodyXml) {
    return '<?xml version="1.0" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="' + $w + '"><w:body>' + $innerBodyXml + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

function Find-ParagraphContaining([object]$doc, [string]$needle) {
    $paras = $doc.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "As a Data Scientist at ELGI Company, I want to ..." paragraph:
#    - "ELGI" becomes its own run containing "X"
#    - "to "/"analyse"/long sentence/" " runs collapse into one run
# ---------------------------------------------------------------------
$p1 = Find-ParagraphContaining $d "Data Scientist at ELGI Company"
$p1InnerXml = (
    '<w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>As</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> a Data Scientist at </w:t></w:r>' +
    '<w:r><w:t>X</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Company, </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>I want</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> to analyse and visualize our sales data for the year 2022, compare it to the actual sales data from 2021, and evaluate the performance of a predictive model we built for customer type classification. </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>So that</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> I get benefits to understand the customer profile in region wise, month wise and product category wise. Helps me to know the staff demand as well.</w:t></w:r>' +
    '</w:p>'
)
$p1.Range.InsertXML((New-Pkg $p1InnerXml))

# ---------------------------------------------------------------------
# 2) "Given a dataset named 'ELGI_2021.csv' ... 'ELGI_2022.csv' ..."
#    -> "ELGI" becomes "Company", each occurrence kept as its own run.
# ---------------------------------------------------------------------
$p2 = Find-ParagraphContaining $d "Given a dataset named 'ELGI_2021.csv'"
$p2InnerXml = (
    '<w:p><w:r><w:t xml:space="preserve">   - Given a dataset named ''</w:t></w:r>' +
    '<w:r><w:t>Company</w:t></w:r>' +
    '<w:r><w:t>_2021.csv'' containing sales data for 2021 and ''</w:t></w:r>' +
    '<w:r><w:t>Company</w:t></w:r>' +
    '<w:r><w:t>_2022.csv'' for 2022,</w:t></w:r>' +
    '</w:p>'
)
$p2.Range.InsertXML((New-Pkg $p2InnerXml))

# ---------------------------------------------------------------------
# 3) "And a trained Random Forest Classifier with 100 estimators," ->
#    expanded sentence listing multiple classifiers, with spell-check
#    proofErr markers around the non-dictionary identifiers.
# ---------------------------------------------------------------------
$p3 = Find-ParagraphContaining $d "And a trained Random Forest Classifier with 100 estimators"
$p3InnerXml = (
    '<w:p>' +
    '<w:r><w:t xml:space="preserve">   - And a trained</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> on multiple classifiers like</w:t></w:r>' +
    '<w:r><w:t>,</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Random Forest (</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>n_estimator</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">=100), Decision tree, SVM (Kernel = linear), K-Nearest </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Neighbors</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>n_neighbors</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> = 5), Logistic Regression, Gradient Boosting (</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>n_estimators</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> = 100)</w:t></w:r>' +
    '</w:p>'
)
$p3.Range.InsertXML((New-Pkg $p3InnerXml))

# ---------------------------------------------------------------------
# 4) Drop the stray <w:lastRenderedPageBreak/> on the "4. Staff
#    Quantity and Cost Analysis:" paragraph.
# ---------------------------------------------------------------------
$p4 = Find-ParagraphContaining $d "Staff Quantity and Cost Analysis"
$p4InnerXml = (
    '<w:p><w:r><w:t xml:space="preserve">4. </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Staff Quantity and Cost Analysis:</w:t></w:r>' +
    '</w:p>'
)
$p4.Range.InsertXML((New-Pkg $p4InnerXml))

Write-Output "done"
